# Apply the benchmark-data update described by the diff:
#  - Rows 2-4 become new MNT-curve entries (MNT159/201/224), each with secparam 512
#    and refreshed timing/byte-size measurements.
#  - The former SS512 rows shift down: the old row 5 (secparam 256) data now
#    appears with secparam 128 in row 5, row 6 (secparam 384) -> 160 in row 6,
#    row 7 (secparam 512) -> 224 in row 7, and three brand-new SS512 rows
#    (256, 384, 512) are appended as rows 8-10.
#  - The sheet's used range grows from A1:R7 to A1:R10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  "MNT159", 512, 20, 20, 20, 0.000782526, 0.005061018, 0.044889295, 0.011978185, 0.040101564, 17697382.4, 17761075.2, 17761075.2, 17761075.2, 17761075.2, 104, 64, 56),
    @(3,  "MNT201", 512, 20, 20, 20, 0.001343608, 0.007150722, 0.057588804, 0.017034423, 0.051809621, 17834188.8, 17864499.2, 17864499.2, 17864499.2, 17864499.2, 104, 64, 56),
    @(4,  "MNT224", 512, 20, 20, 20, 0.001549947, 0.008833385,  0.070989013, 0.019518232, 0.062476432, 17886822.4, 17905664,   17905664,   17905664,   17905664,   104, 64, 56),
    @(5,  "SS512",  128, 20, 20, 20, 0.001719761, 0.01428448,   0.014179206, 0.016858089, 0.012849689, 17989632,   17989632,   17989632,   17989632,   17989632,   104, 64, 56),
    @(6,  "SS512",  160, 20, 20, 20, 0.001662135, 0.012737262,  0.013528812, 0.01528455,  0.011559737, 17989632,   17989632,   17989632,   17989632,   17989632,   104, 64, 56),
    @(7,  "SS512",  224, 20, 20, 20, 0.001555836, 0.012916672,  0.013057399, 0.015018916, 0.011732829, 17989632,   17989632,   17989632,   17989632,   17989632,   104, 64, 56),
    @(8,  "SS512",  256, 20, 20, 20, 0.001542807, 0.014727831,  0.014608085, 0.016889131, 0.01244874,  17989632,   17989632,   17989632,   17989632,   17989632,   104, 64, 56),
    @(9,  "SS512",  384, 20, 20, 20, 0.001709652, 0.014466238,  0.014772856, 0.016612387, 0.012647295, 17989632,   17989632,   17989632,   17989632,   17989632,   104, 64, 56),
    @(10, "SS512",  512, 20, 20, 20, 0.001668262, 0.014079392,  0.013857687, 0.016103268, 0.012204862, 17989632,   17989632,   17989632,   17989632,   17989632,   104, 64, 56)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($col = 1; $col -lt $row.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
